$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$lastCol = $ws.Cells.Item(1, $ws.Columns.Count).End(-4159).Column

for ($r = 2; $r -le $lastRow; $r++) {
    for ($c = 2; $c -le $lastCol; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value2
        if ($val -ne $null) {
            $cell.Value = [math]::Floor([double]$val + 0.5)
        }
    }
}
